# Add a new data row to the "Removal" sheet (copy of the existing last row,
# row 5, with the testcase name swapped out) and leave that sheet as the
# active one with C8 selected -- mirroring how the source workbook was
# edited by hand in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Removal")

# Duplicate row 5 (keeps every cell value + style) into a brand new row 6.
$ws.Rows(5).Copy() | Out-Null
$ws.Rows(6).Insert() | Out-Null

# Row 6 is a new, distinct test case -- only the first column differs.
$ws.Range("A6").Value = "testT4258"

# Make "Removal" the active sheet/tab, with C8 selected, matching the
# workbook's saved view state.
$ws.Activate() | Out-Null
$ws.Range("C8").Select() | Out-Null
